$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Mechanical Equipment"
$ws.Range("B4").Value = "TEST_serienummer"
$ws.Range("C4").Value = "Serienummer"

$ws.Range("G8").Select() | Out-Null
